$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number + report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Cells switching FROM numeric TO the "no data" text placeholder ---
# (copy format+text from the row-14 header template cell in the same column,
#  which already holds the correct shared-string placeholder + style)
$ws.Range("C14").Copy($ws.Range("C17"))
$ws.Range("D14").Copy($ws.Range("D17"))
$ws.Range("E14").Copy($ws.Range("E17"))
$ws.Range("F14").Copy($ws.Range("F22"))
$ws.Range("G14").Copy($ws.Range("G30"))
$ws.Range("H14").Copy($ws.Range("H30"))

# --- Cells switching FROM the text placeholder TO numeric data ---
# (copy number format from a same-column numeric sibling cell, then set the value)
$ws.Range("D16").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("E16").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("D16").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("E16").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -28.571428571428
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 13.333333333333
$ws.Range("L16").Value = 70
$ws.Range("M16").Value = 30.769230769230
$ws.Range("N16").Value = -84.451219512195
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = -12.5
$ws.Range("L17").Value = -14.754098360655
$ws.Range("M17").Value = 44.444444444444
$ws.Range("N17").Value = -46.938775510204
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = -45.614035087719
$ws.Range("L18").Value = -7.462686567164
$ws.Range("M18").Value = 6.896551724137
$ws.Range("N18").Value = -89.869281045751
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 52
$ws.Range("H19").Value = 23.809523809523
$ws.Range("I19").Value = 399
$ws.Range("J19").Value = 341
$ws.Range("K19").Value = 17.008797653958
$ws.Range("L19").Value = 55.252918287937
$ws.Range("M19").Value = 5.835543766578
$ws.Range("N19").Value = -65.984654731457
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = -46.153846153846
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = -3.125
$ws.Range("L20").Value = 3.333333333333
$ws.Range("M20").Value = 121.428571428571
$ws.Range("N20").Value = -91.014492753623
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -12
$ws.Range("F21").Value = 80
$ws.Range("H21").Value = -14.893617021276
$ws.Range("I21").Value = 600
$ws.Range("J21").Value = 597
$ws.Range("K21").Value = 0.502512562814
$ws.Range("L21").Value = 34.228187919463
$ws.Range("M21").Value = 13.851992409867
$ws.Range("N21").Value = -76.589933671478
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 14.285714285714
$ws.Range("M22").Value = -15.789473684210
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 140
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 12.903225806451
$ws.Range("I24").Value = 576
$ws.Range("J24").Value = 683
$ws.Range("K24").Value = -15.666178623718
$ws.Range("L24").Value = 6.077348066298
$ws.Range("M24").Value = 69.411764705882
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 30.769230769230
$ws.Range("I25").Value = 143
$ws.Range("J25").Value = 125
$ws.Range("K25").Value = 14.4
$ws.Range("L25").Value = 74.390243902439
$ws.Range("M25").Value = 9.160305343511
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -22.222222222222
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 32
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = -11.111111111111
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = -33.333333333333
$ws.Range("L30").Value = 0
